$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AJ2").Value = 970
$ws.Range("O2").Value = 1.27
$ws.Range("F3").Value = 1.85
$ws.Range("G3").Value = 2.06
$ws.Range("H3").Value = 4.5
$ws.Range("K3").Value = 3.8
$ws.Range("N3").Value = 2.98
$ws.Range("O3").Value = 1.42
$ws.Range("S3").Value = 3.7
$ws.Range("T3").Value = 1.95
$ws.Range("U3").Value = 1.83
$ws.Range("W3").Value = 1.94
$ws.Range("AF4").Value = 12.5
$ws.Range("AJ4").Value = 23
$ws.Range("AK4").Value = 22
$ws.Range("AN4").Value = 14.5
$ws.Range("P4").Value = 1.93
$ws.Range("AC5").Value = 9.6
$ws.Range("AE5").Value = 48
$ws.Range("AF5").Value = 19.5
$ws.Range("AH5").Value = 21
$ws.Range("AJ5").Value = 42
$ws.Range("AK5").Value = 34
$ws.Range("AL5").Value = 50
$ws.Range("AO5").Value = 46
$ws.Range("G5").Value = 2.62
$ws.Range("I5").Value = 3.55
$ws.Range("J5").Value = 3.3
$ws.Range("K5").Value = 3.75
$ws.Range("L5").Value = 1.34
$ws.Range("N5").Value = 3.5
$ws.Range("Q5").Value = 1.95
$ws.Range("R5").Value = 1.33
$ws.Range("V5").Value = 1.39
$ws.Range("X5").Value = 17
$ws.Range("G6").Value = 1.41
$ws.Range("L6").Value = 1.25
$ws.Range("R6").Value = 1.7
$ws.Range("S6").Value = 2.12
$ws.Range("T6").Value = 1.78
$ws.Range("U6").Value = 2.06
$ws.Range("AA7").Value = 580
$ws.Range("AB7").Value = 10.5
$ws.Range("AD7").Value = 46
$ws.Range("AG7").Value = 11
$ws.Range("AK7").Value = 16.5
$ws.Range("AL7").Value = 36
$ws.Range("F7").Value = 1.27
$ws.Range("K7").Value = 7.4
$ws.Range("P7").Value = 2.68
$ws.Range("R7").Value = 1.67
$ws.Range("S7").Value = 2.32
$ws.Range("T7").Value = 2.02
$ws.Range("U7").Value = 1.86
$ws.Range("X7").Value = 1000
$ws.Range("Y7").Value = 55
$ws.Range("N8").Value = 3.35
$ws.Range("P8").Value = 1.76
$ws.Range("R8").Value = 1.28
$ws.Range("AH9").Value = 21
$ws.Range("T9").Value = 1.69
$ws.Range("T10").Value = 1.9
$ws.Range("AF11").Value = 46
$ws.Range("AH11").Value = 23
$ws.Range("AI11").Value = 40
$ws.Range("G11").Value = 6
$ws.Range("H11").Value = 1.7
$ws.Range("I11").Value = 1.77
$ws.Range("K11").Value = 4.2
$ws.Range("P11").Value = 1.84
$ws.Range("V11").Value = 2.28
$ws.Range("AA12").Value = 14
$ws.Range("AB12").Value = 48
$ws.Range("AE12").Value = 14.5
$ws.Range("F12").Value = 8.800000000000001
$ws.Range("G12").Value = 11
$ws.Range("H12").Value = 1.31
$ws.Range("I12").Value = 1.41
$ws.Range("J12").Value = 5.6
$ws.Range("K12").Value = 7
$ws.Range("N12").Value = 5.7
$ws.Range("P12").Value = 2.62
$ws.Range("Q12").Value = 1.5
$ws.Range("R12").Value = 1.64
$ws.Range("S12").Value = 2.32
$ws.Range("T12").Value = 1.7
$ws.Range("U12").Value = 2.02
$ws.Range("V12").Value = 3.45
$ws.Range("W12").Value = 1.1
$ws.Range("Z12").Value = 9.6
$ws.Range("U13").Value = 1.61
$ws.Range("X13").Value = 10.5
$ws.Range("N14").Value = 1.1
$ws.Range("P14").Value = 1.25
$ws.Range("R14").Value = 1.18
$ws.Range("T14").Value = 1.04
$ws.Range("U14").Value = 1.04
$ws.Range("AA15").Value = 0
$ws.Range("AB15").Value = 0
$ws.Range("AC15").Value = 0
$ws.Range("AD15").Value = 0
$ws.Range("AE15").Value = 0
$ws.Range("AF15").Value = 0
$ws.Range("AG15").Value = 0
$ws.Range("AH15").Value = 0
$ws.Range("AI15").Value = 0
$ws.Range("AJ15").Value = 0
$ws.Range("AK15").Value = 0
$ws.Range("AL15").Value = 0
$ws.Range("AM15").Value = 0
$ws.Range("AN15").Value = 0
$ws.Range("AO15").Value = 0
$ws.Range("D15").Value = "Crotone"
$ws.Range("E15").Value = "S.S.D. Casarano Calcio"
$ws.Range("F15").Value = 1.58
$ws.Range("G15").Value = 1.75
$ws.Range("H15").Value = 6
$ws.Range("I15").Value = 7.6
$ws.Range("J15").Value = 3.4
$ws.Range("K15").Value = 4.4
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 3.55
$ws.Range("O15").Value = 1.29
$ws.Range("P15").Value = 1.86
$ws.Range("Q15").Value = 1.89
$ws.Range("R15").Value = 0
$ws.Range("S15").Value = 0
$ws.Range("T15").Value = 0
$ws.Range("U15").Value = 0
$ws.Range("V15").Value = 0
$ws.Range("W15").Value = 0
$ws.Range("X15").Value = 0
$ws.Range("Y15").Value = 0
$ws.Range("Z15").Value = 0
$ws.Range("AA16").Value = 34
$ws.Range("AB16").Value = 12
$ws.Range("AC16").Value = 7.8
$ws.Range("AD16").Value = 12
$ws.Range("AE16").Value = 32
$ws.Range("AF16").Value = 27
$ws.Range("AG16").Value = 16
$ws.Range("AH16").Value = 22
$ws.Range("AI16").Value = 55
$ws.Range("AJ16").Value = 80
$ws.Range("AK16").Value = 55
$ws.Range("AL16").Value = 70
$ws.Range("AM16").Value = 1000
$ws.Range("AN16").Value = 1000
$ws.Range("AO16").Value = 1000
$ws.Range("D16").Value = "Pesaro"
$ws.Range("E16").Value = "Ravenna"
$ws.Range("F16").Value = 3.4
$ws.Range("G16").Value = 3.95
$ws.Range("H16").Value = 2.26
$ws.Range("I16").Value = 2.52
$ws.Range("J16").Value = 3.05
$ws.Range("K16").Value = 3.75
$ws.Range("L16").Value = 1.5
$ws.Range("M16").Value = 1.1
$ws.Range("N16").Value = 2.92
$ws.Range("O16").Value = 1.42
$ws.Range("P16").Value = 1.66
$ws.Range("Q16").Value = 2.22
$ws.Range("R16").Value = 1.24
$ws.Range("S16").Value = 4.3
$ws.Range("T16").Value = 1.93
$ws.Range("U16").Value = 1.89
$ws.Range("V16").Value = 1.66
$ws.Range("W16").Value = 1.34
$ws.Range("X16").Value = 12
$ws.Range("Y16").Value = 8.6
$ws.Range("Z16").Value = 15
$ws.Range("G17").Value = 4.6
$ws.Range("I17").Value = 2.26
$ws.Range("J17").Value = 3.2
$ws.Range("S17").Value = 4.3
$ws.Range("T17").Value = 1.88
$ws.Range("U17").Value = 1.67
$ws.Range("V17").Value = 1.8
$ws.Range("X17").Value = 10
$ws.Range("AN18").Value = 1000
$ws.Range("K18").Value = 5.2
$ws.Range("P18").Value = 1.76
$ws.Range("Q18").Value = 1.96
$ws.Range("R19").Value = 1.48
$ws.Range("U19").Value = 1.93
$ws.Range("Z19").Value = 75
$ws.Range("L20").Value = 1.47
$ws.Range("R20").Value = 1.25
$ws.Range("T20").Value = 1.94
$ws.Range("N21").Value = 3.15
$ws.Range("P21").Value = 1.72
$ws.Range("Q21").Value = 2.34
$ws.Range("AE22").Value = 44
$ws.Range("AN22").Value = 12.5
$ws.Range("G22").Value = 1.94
$ws.Range("H22").Value = 4.2
$ws.Range("I22").Value = 4.5
$ws.Range("J22").Value = 3.85
$ws.Range("N22").Value = 4.1
$ws.Range("P22").Value = 2.06
$ws.Range("R22").Value = 1.41
$ws.Range("S22").Value = 3
$ws.Range("V22").Value = 1.29
$ws.Range("W22").Value = 2.06
$ws.Range("X22").Value = 16.5
$ws.Range("Y22").Value = 21
$ws.Range("AC23").Value = 10
$ws.Range("AD23").Value = 16.5
$ws.Range("AK23").Value = 16
$ws.Range("AN23").Value = 8.4
$ws.Range("AO23").Value = 32
$ws.Range("J23").Value = 4.3
$ws.Range("K23").Value = 4.4
$ws.Range("P23").Value = 2.58
$ws.Range("Q23").Value = 1.62
$ws.Range("R23").Value = 1.63
$ws.Range("G24").Value = 2.78
$ws.Range("Q24").Value = 2.02
$ws.Range("AE25").Value = 550
$ws.Range("J25").Value = 8.199999999999999
$ws.Range("K25").Value = 8.800000000000001
$ws.Range("M25").Value = 1.03
$ws.Range("Z25").Value = 260
$ws.Range("AF26").Value = 30
$ws.Range("AH26").Value = 16.5
$ws.Range("AM26").Value = 65
$ws.Range("AO26").Value = 13.5
$ws.Range("G26").Value = 3.6
$ws.Range("I26").Value = 2.18
$ws.Range("J26").Value = 4
$ws.Range("K26").Value = 4.5
$ws.Range("M26").Value = 1.01
$ws.Range("N26").Value = 5.3
$ws.Range("O26").Value = 1.18
$ws.Range("P26").Value = 1.62
$ws.Range("Q26").Value = 1.56
$ws.Range("R26").Value = 1.48
$ws.Range("S26").Value = 2.08
$ws.Range("T26").Value = 1.04
$ws.Range("U26").Value = 2.5
$ws.Range("V26").Value = 1.84
$ws.Range("X26").Value = 24
$ws.Range("Y26").Value = 15.5
